$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.607.72'
$ws.Range("E2").Value = '  -0.35%  '

$ws.Range("D3").Value = '2.300.85'
$ws.Range("E3").Value = '  -0.13%  '

$ws.Range("E4").Value = '  +0.18%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '310.30'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -2.01%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '104.92'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.01%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.626'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.53%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").Value = '  -0.45%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '39.62'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.54%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0910'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.40%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '8.33'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.20%  '

$ws.Range("E13").Value = '  +0.01%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.987'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.64%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '15.24'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.32%  '

$ws.Range("D16").Value = '2.651.39'
$ws.Range("E16").Value = '  -0.03%  '

$ws.Range("D17").Value = '2.301.38'
$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("D18").Value = '42.827.12'
$ws.Range("E18").Value = '  +0.38%  '

$ws.Range("E19").Value = '  -3.19%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.77'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.59%  '

$ws.Range("E21").Value = '  -1.09%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '73.43'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.92%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '3.45'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.71%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '269.16'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.27%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.24'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.07%  '

$ws.Range("B26").Value = 'Filecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '8.02'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +19.63%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.02%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '10.94'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("E29").Value = '  -2.06%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '37.92'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.52%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '22.19'
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '165.79'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.04%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0863'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.78%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.79'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +2.10%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.131'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.87%  '

$ws.Range("E36").Value = '  -0.80%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.63'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.12%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.0357'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.64%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.79'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.53%  '

$ws.Range("E40").Value = '  -4.22%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '107.74'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +12.55%  '

$ws.Range("E42").Value = '  -4.33%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '71.08'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.32%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.229'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +0.48%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '12.21'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -3.13%  '

$ws.Range("D47").Value = '1.693.63'
$ws.Range("E47").Value = '  +1.58%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '110.76'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -5.52%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '75.60'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -5.87%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '8.86'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.42%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '5.17'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -2.41%  '
